# Swap the species-observation data between row 5 and row 6, while keeping
# the shared/location fields (C, D, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AW, AX, AY) unchanged on each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "J", "M", "Q", "R", "AF")

foreach ($col in $cols) {
    $addr5 = "$col" + "5"
    $addr6 = "$col" + "6"
    $v5 = $ws.Range($addr5).Value2
    $v6 = $ws.Range($addr6).Value2
    $ws.Range($addr5).Value = $v6
    $ws.Range($addr6).Value = $v5
}
